$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty "street" column (E) -- shifts price/unit_price/land_size/garage left
$ws.Columns("E").Delete()

# Header row
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "city"
$ws.Range("C1").Value = "type"
$ws.Range("D1").Value = "size"
$ws.Range("E1").Value = "price"
$ws.Range("F1").Value = "unit_price"
$ws.Range("G1").Value = "land_size"
$ws.Range("H1").Value = "count"

# Row 2
$ws.Range("A2").Value = 43841
$ws.Range("B2").Value = "Algyo"
$ws.Range("C2").Value = "House"
$ws.Range("D2").Value = 118.882352941177
$ws.Range("E2").Value = 30002941.1764706
$ws.Range("F2").Value = 278100.389936275
$ws.Range("H2").Value = 18

# Row 3
$ws.Range("A3").Value = 43841
$ws.Range("B3").Value = "Morahalom"
$ws.Range("C3").Value = "House"
$ws.Range("D3").Value = 118.179487179487
$ws.Range("E3").Value = 24864102.5641026
$ws.Range("F3").Value = 215486.314278871
$ws.Range("H3").Value = 40

# Row 4
$ws.Range("A4").Value = 43841
$ws.Range("B4").Value = "Szeged"
$ws.Range("C4").Value = "Garage"
$ws.Range("D4").Value = 17.7551020408163
$ws.Range("E4").Value = 4165918.36734694
$ws.Range("F4").Value = 242099.910224355
$ws.Range("H4").Value = 99

# Row 5
$ws.Range("A5").Value = 43841
$ws.Range("B5").Value = "Szeged"
$ws.Range("C5").Value = "House"
$ws.Range("D5").Value = 160.991161616162
$ws.Range("E5").Value = 54929393.9393939
$ws.Range("F5").Value = 992379.652400437
$ws.Range("G5").Value = 470.02398989899
$ws.Range("H5").Value = 770

# Row 6
$ws.Range("A6").Value = 43843
$ws.Range("B6").Value = "Algyo"
$ws.Range("C6").Value = "House"
$ws.Range("D6").Value = 118.944444444444
$ws.Range("E6").Value = 29663888.8888889
$ws.Range("F6").Value = 273715.183087963
$ws.Range("G6").Value = 513.833333333333
$ws.Range("H6").Value = 18

# Row 7
$ws.Range("A7").Value = 43843
$ws.Range("B7").Value = "Morahalom"
$ws.Range("C7").Value = "House"
$ws.Range("D7").Value = 117.475
$ws.Range("E7").Value = 24720000
$ws.Range("F7").Value = 215462.955705053
$ws.Range("G7").Value = 157.05
$ws.Range("H7").Value = 40

# Row 8
$ws.Range("A8").Value = 43843
$ws.Range("B8").Value = "Szeged"
$ws.Range("C8").Value = "Garage"
$ws.Range("D8").Value = 17.6969696969697
$ws.Range("E8").Value = 4154141.41414141
$ws.Range("F8").Value = 242309.862182307
$ws.Range("H8").Value = 99

# Row 9
$ws.Range("A9").Value = 43844
$ws.Range("B9").Value = "Algyo"
$ws.Range("C9").Value = "House"
$ws.Range("D9").Value = 118.9444444444444
$ws.Range("E9").Value = 29663888.88888889
$ws.Range("F9").Value = 273715.1830879633
$ws.Range("G9").Value = 513.8333333333334
$ws.Range("H9").Value = 18

# Row 10
$ws.Range("A10").Value = 43844
$ws.Range("B10").Value = "Morahalom"
$ws.Range("C10").Value = "House"
$ws.Range("D10").Value = 117.475
$ws.Range("E10").Value = 24720000
$ws.Range("F10").Value = 215462.9557050533
$ws.Range("G10").Value = 157.05
$ws.Range("H10").Value = 40

# Row 11
$ws.Range("A11").Value = 43844
$ws.Range("B11").Value = "Szeged"
$ws.Range("C11").Value = "Garage"
$ws.Range("D11").Value = 17.6969696969697
$ws.Range("E11").Value = 4154141.414141414
$ws.Range("F11").Value = 242309.862182307
$ws.Range("H11").Value = 99

# Match date-cell formatting (border/numfmt/alignment/font) for the newly added rows
$ws.Range("A2").Copy()
$ws.Range("A6:A11").PasteSpecial(-4122)
